$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F; this shifts old F (District) to G
$ws.Columns("F").Insert()

# Set header for new Address column
$ws.Range("F2").Value = "Address"

# Populate Address values extracted from the school/location text in column B
$ws.Range("F3").Value = "G U H S AtnoorAfzalpur"
$ws.Range("F4").Value = "G P U C (H S) Mugthihalli"
$ws.Range("F5").Value = "S M R High School HebbalagereChannagiri"
$ws.Range("F9").Value = "Govt. Adarsha Vidyalaya Chittapur"
$ws.Range("F10").Value = "G H S KodliChincholi"
$ws.Range("F11").Value = "G J C H S PanchanahalliKadur"
$ws.Range("F12").Value = "G H S HarkudBasavakalyan"
$ws.Range("F13").Value = "G H S ManthalBadavakalyan"
$ws.Range("F14").Value = "G H S JeganehalliBirur Block Kadur"
$ws.Range("F15").Value = "S B D H S Arsikere"
$ws.Range("F16").Value = "Nalanda P U CollegeHigh School SectionSagar"
$ws.Range("F17").Value = "Mothiveerappa Govt. P U College"
$ws.Range("F18").Value = "Govt. High School KalgiChittapur"
$ws.Range("F19").Value = "S V P S B C Girls High School Badami"
$ws.Range("F22").Value = "G H S Ekamba Aurad"
$ws.Range("F23").Value = "S S High School AgarkhedIndi"
$ws.Range("F24").Value = "S M S High School YadavanahallyArsikere"
$ws.Range("F25").Value = "Govt. High School Veerapur"
$ws.Range("F26").Value = "G J C (H S) Yagati Kadur"
$ws.Range("F27").Value = "Vishvodaya High School Krishnaiahna DoddiKanakapura"
$ws.Range("F28").Value = "Shree J B K High School Basavakalyan"
$ws.Range("F29").Value = "Govt. H S MurkiAurad"
$ws.Range("F30").Value = "R M B P U College (H S) AthargaIndi"
$ws.Range("F31").Value = "Malanad High School GoutampurSagar"
$ws.Range("F32").Value = "Sri Kalmarudeswara High SchoolMarle"
$ws.Range("F33").Value = "CholachaguddBadami"
$ws.Range("F34").Value = "G H S GadikeshwarChincholi"
$ws.Range("F35").Value = "G B H S ShahabadChittapur"
$ws.Range("F36").Value = "Matru Manidr H S Chincholi"
$ws.Range("F37").Value = "Bharthi Viswaseva SadanHigh School Somapur"
$ws.Range("F38").Value = "G H S NimbargaAland"
$ws.Range("F39").Value = "S B P U College BolegaonIndi"
$ws.Range("F40").Value = "G H S Satanoor Chittapur"
$ws.Range("F41").Value = "V S High School KedalasaraSagar"
$ws.Range("F42").Value = "G P U C HalebeeduBelur"
$ws.Range("F43").Value = "S M High School TerdalJamakhandi"
$ws.Range("F45").Value = "HasadurgaKanakapura"
$ws.Range("F46").Value = "G J C JodihochihalliBeerur BlockKadur"
$ws.Range("F47").Value = "S S P U College (H S) Indi"
$ws.Range("F48").Value = "G J C Kadur"
$ws.Range("F49").Value = "SchoolHubli"
